# Applies the commit "Update gh-pages to output generated at 456a3b4" to
# the 杭州-漫展信息 workbook.
#
# Summary of the change:
#  - Sheet "展览" (exhibition list): the event row for 2024-05-04
#    "杭州·Ani idol08偶像剧场" (row 4) was removed entirely, which shifts
#    every following row up by one and drops the former last row
#    (2024-09-15 item, formerly row 42). The "想去人数" (F) counters were
#    also refreshed (bumped) for the surviving rows.
#  - Sheet "演出" (performances): the "想去人数" counter for the
#    2024-06-22 "时光代理人" show (row 15) was refreshed.
#  - Sheet "全部类型" (all types): the "想去人数" counters for the rows
#    that mirror the surviving 展览 events were refreshed to the same
#    new values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) 展览 sheet: delete the obsolete row (2024-05-04 Ani idol08) and
#    refresh the "想去人数" (F) values for the rows that shifted up.
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")

# Deleting the entire row shifts rows 5..42 up to 4..41 and keeps every
# other cell (including the literal serial numbers in column A) intact,
# matching how Excel performs a native row delete.
$wsExpo.Rows.Item(4).Delete()

$expoUpdates = @(
    @(4, 97),
    @(5, 102),
    @(6, 280),
    @(8, 1138),
    @(9, 396),
    @(10, 89),
    @(11, 97),
    @(12, 123),
    @(13, 21),
    @(14, 247),
    @(15, 142),
    @(16, 133),
    @(17, 1291),
    @(18, 490),
    @(19, 180),
    @(20, 293),
    @(22, 656),
    @(23, 1062),
    @(25, 1941),
    @(26, 2456),
    @(27, 1222),
    @(28, 51),
    @(29, 166),
    @(30, 345),
    @(31, 705),
    @(32, 749),
    @(33, 883),
    @(34, 111),
    @(36, 739),
    @(37, 316),
    @(38, 569),
    @(39, 688),
    @(40, 301),
    @(41, 197)
)

foreach ($pair in $expoUpdates) {
    $row = $pair[0]
    $val = $pair[1]
    $wsExpo.Range("F" + $row).Value2 = $val
}

# ---------------------------------------------------------------------
# 2) 演出 sheet: refresh the "想去人数" value for row 15.
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F15").Value2 = 294

# ---------------------------------------------------------------------
# 3) 全部类型 sheet: refresh "想去人数" values for the rows that mirror
#    the surviving 展览 events.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$allUpdates = @(
    @(3, 97),
    @(7, 102),
    @(8, 280),
    @(12, 1138),
    @(13, 396),
    @(14, 89),
    @(15, 97),
    @(16, 123),
    @(17, 247),
    @(19, 142),
    @(20, 133),
    @(21, 1291),
    @(22, 490),
    @(23, 180),
    @(24, 293),
    @(26, 1062),
    @(27, 2456),
    @(29, 1222),
    @(30, 51),
    @(34, 166),
    @(35, 345),
    @(36, 705),
    @(39, 749),
    @(40, 883),
    @(41, 739),
    @(42, 316),
    @(43, 569),
    @(44, 688),
    @(45, 301),
    @(48, 197)
)

foreach ($pair in $allUpdates) {
    $row = $pair[0]
    $val = $pair[1]
    $wsAll.Range("F" + $row).Value2 = $val
}
